$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "full random"
$ws.Range("B5").Value = "sentences"
$ws.Range("C5").Value = 3000
$ws.Range("D5").Value = 200
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = "2000s (ca)"
$ws.Range("G5").Value = "auto"
$ws.Range("H5").Value = 8
$ws.Range("I5").Value = "bad"

$ws.Range("I10").Select()
